$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order (Oyuncu Adı, Pozisyon, Takım) for rows 2..18
$rows = @(
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Nick Richards", "C", "Phoenix Suns"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
